$wb = $excel.ActiveWorkbook

$txSheets = @(
    "O_TransactionActivity",
    "O_TransactionActivity2",
    "O_TransactionActivity3",
    "O_TransactionActivity5",
    "O_TransactionActivity6"
)

foreach ($name in $txSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("D2").Value = "Source ASC"
}
